$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Insert two new rows for "bird" and "goodegg" right after the "sputnick" row (27) ---
# This pushes the "\YFCF" block (old rows 30,32,33,34) down to rows 32,34,35,36.
$ws.Rows("28:29").Insert()

# New row 28: bird
$ws.Range("A28").Value2 = "bird"
$ws.Range("B28").Value2 = "https://web.archive.org/web/19970719105040im_/http://dewey.rug.ac.be/barn/tex/bird.html"

# New row 29: goodegg
$ws.Range("A29").Value2 = "goodegg"
$ws.Range("B29").Value2 = "https://web.archive.org/web/19970719105100im_/http://dewey.rug.ac.be/barn/tex/goodegg.html"

# --- Update the "video" row URL (row 7) to the new archive snapshot ---
$ws.Range("B7").Value2 = "https://web.archive.org/web/19970715175954im_/http://dewey.rug.ac.be/barn/tex/video.html"

# --- Re-create the pre-existing HotSonic hyperlink, whose anchor needs to move from B33 to B35 ---
# (Target URL kept byte-for-byte identical to the original relationship.)
$ws.Range("B33").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B35"), "https://web.archive.org/web/19970719105949/http:/dewey.rug.ac.be/YFCF/HotSonic.html") | Out-Null

# --- Add the two new hyperlinks (yfcf row 5, video row 7) ---
$ws.Hyperlinks.Add($ws.Range("B5"), "https://web.archive.org/web/19990220153326im_/http://dewey.rug.ac.be/barn/tex/yfcf.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://web.archive.org/web/19970715175954im_/http://dewey.rug.ac.be/barn/tex/video.html") | Out-Null

# Restore the plain "Hyperlink" cell style (rather than the ad-hoc style Hyperlinks.Add synthesizes)
$ws.Range("B35").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"

# --- Selection moves to B32 in the final sheet ---
$ws.Range("B32").Select()
